$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: FERNANDO LUIS BUELVAS RODRIGUEZ, period 2201, 26666 / 1000000
$ws.Range("C16").Value = "1143341397"
$ws.Range("D16").Value = "FERNANDO LUIS BUELVAS RODRIGUEZ"
$ws.Range("E16").Value = "2201"
$ws.Range("F16").Value = 26666
$ws.Range("G16").Value = 1000000

# Row 17: MARIO MIGUEL SALAS COAVAS, period 2201, 64000 / 2000000
$ws.Range("C17").Value = "1143334296"
$ws.Range("D17").Value = "MARIO MIGUEL SALAS COAVAS"
$ws.Range("E17").Value = "2201"
$ws.Range("F17").Value = 64000
$ws.Range("G17").Value = 2000000

# Row 18: FERNANDO LUIS BUELVAS RODRIGUEZ, period 2202, 40000 / 1000000
$ws.Range("C18").Value = "1143341397"
$ws.Range("D18").Value = "FERNANDO LUIS BUELVAS RODRIGUEZ"
$ws.Range("E18").Value = "2202"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

# Row 19: FERNANDO LUIS BUELVAS RODRIGUEZ, period 2203, 40000 / 1000000
$ws.Range("C19").Value = "1143341397"
$ws.Range("D19").Value = "FERNANDO LUIS BUELVAS RODRIGUEZ"
$ws.Range("E19").Value = "2203"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

# Row 20: LUZ KARIME STEVENSON BARCHA, period 2207, 153760 / 3844000 (only G changes)
$ws.Range("G20").Value = 3844000
